$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 204: new "Source: 2021-02-30" marker row (only column B populated)
$ws.Range("B204").Value = "Source: 2021-02-30"

# Row 205-210: week 10 data
$ws.Range("A205").Value = 2021
$ws.Range("B205").Value = 10
$ws.Range("C205").Value = "0-4"
$ws.Range("D205").Value = 18750
$ws.Range("E205").Value = 5.1

$ws.Range("A206").Value = 2021
$ws.Range("B206").Value = 10
$ws.Range("C206").Value = "5-14"
$ws.Range("D206").Value = 28009
$ws.Range("E206").Value = 8.1

$ws.Range("A207").Value = 2021
$ws.Range("B207").Value = 10
$ws.Range("C207").Value = "15-34"
$ws.Range("D207").Value = 108796
$ws.Range("E207").Value = 7.2

$ws.Range("A208").Value = 2021
$ws.Range("B208").Value = 10
$ws.Range("C208").Value = "35-59"
$ws.Range("D208").Value = 161342
$ws.Range("E208").Value = 7

$ws.Range("A209").Value = 2021
$ws.Range("B209").Value = 10
$ws.Range("C209").Value = "60-79"
$ws.Range("D209").Value = 89815
$ws.Range("E209").Value = 5.7

$ws.Range("A210").Value = 2021
$ws.Range("B210").Value = 10
$ws.Range("C210").Value = ">=80"
$ws.Range("D210").Value = 44444
$ws.Range("E210").Value = 5

# Row 211 intentionally left blank (gap in source data)

# Row 212-217: week 11 data
$ws.Range("A212").Value = 2021
$ws.Range("B212").Value = 11
$ws.Range("C212").Value = "0-4"
$ws.Range("D212").Value = 28472
$ws.Range("E212").Value = 5.1

$ws.Range("A213").Value = 2021
$ws.Range("B213").Value = 11
$ws.Range("C213").Value = "5-14"
$ws.Range("D213").Value = 41666
$ws.Range("E213").Value = 8

$ws.Range("A214").Value = 2021
$ws.Range("B214").Value = 11
$ws.Range("C214").Value = "15-34"
$ws.Range("D214").Value = 119444
$ws.Range("E214").Value = 8.6

$ws.Range("A215").Value = 2021
$ws.Range("B215").Value = 11
$ws.Range("C215").Value = "35-59"
$ws.Range("D215").Value = 170601
$ws.Range("E215").Value = 8.5

$ws.Range("A216").Value = 2021
$ws.Range("B216").Value = 11
$ws.Range("C216").Value = "60-79"
$ws.Range("D216").Value = 91890
$ws.Range("E216").Value = 6.6

$ws.Range("A217").Value = 2021
$ws.Range("B217").Value = 11
$ws.Range("C217").Value = ">=80"
$ws.Range("D217").Value = 44444
$ws.Range("E217").Value = 5.2

# Update frozen-pane / selection view state to match new data extent
$ws.Application.ActiveWindow.ScrollRow = 194
$ws.Range("E215").Select() | Out-Null
